$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.807599666666667
$ws.Range("H2").Value = 5.422799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 170.93328
$ws.Range("N2").Value = 512.79984
$ws.Range("O2").Value = 0.7687311215213114
$ws.Range("P2").Value = 0.7687311215213115
$ws.Range("Q2").Value = 308.97893995024
$ws.Range("R2").Value = 2780.81045955216
$ws.Range("S2").Value = 0.7687311215213114
$ws.Range("T2").Value = 0.7687311215213115

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.807599666666667
$ws.Range("H3").Value = 5.422799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.31217066666667
$ws.Range("N3").Value = 120.936512
$ws.Range("O3").Value = 0.1812942463137967
$ws.Range("P3").Value = 0.1812942463137967
$ws.Range("Q3").Value = 72.86826625967645
$ws.Range("R3").Value = 655.8143963370881
$ws.Range("S3").Value = 0.1812942463137967
$ws.Range("T3").Value = 0.1812942463137967

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.807599666666667
$ws.Range("H4").Value = 5.422799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 11.112244
$ws.Range("N4").Value = 33.336732
$ws.Range("O4").Value = 0.04997463216489184
$ws.Range("P4").Value = 0.04997463216489184
$ws.Range("Q4").Value = 20.08648855031867
$ws.Range("R4").Value = 180.778396952868
$ws.Range("S4").Value = 0.04997463216489184
$ws.Range("T4").Value = 0.04997463216489184

